# Rewrites the BOQ line-item rows (8-16) and the grand-total cells (G18/H18/G20/H20)
# so the sheet matches the target revision of the "3rdRunningNoExtra_iter2" bill.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8,1).Value = "'"
$ws.Cells.Item(8,3).Value = 37
$ws.Cells.Item(8,4).Value = "'1.0"
$ws.Cells.Item(8,5).Value = "'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it's ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = "'0.00"

# Row 9
$ws.Cells.Item(9,1).Value = "'"
$ws.Cells.Item(9,3).Value = 78
$ws.Cells.Item(9,4).Value = "'2.0"
$ws.Cells.Item(9,5).Value = "'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = "'0.00"

# Row 10
$ws.Cells.Item(10,1).Value = "'P. point"
$ws.Cells.Item(10,3).Value = 91
$ws.Cells.Item(10,4).Value = "'6"
$ws.Cells.Item(10,5).Value = "'On board"
$ws.Cells.Item(10,6).Value = 136
$ws.Cells.Item(10,7).Value = "'12376.00"

# Row 11
$ws.Cells.Item(11,1).Value = "'"
$ws.Cells.Item(11,3).Value = 78
$ws.Cells.Item(11,4).Value = "'4.0"
$ws.Cells.Item(11,5).Value = "'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Cells.Item(11,6).Value = 50
$ws.Cells.Item(11,7).Value = "'3900.00"

# Row 12
$ws.Cells.Item(12,1).Value = "'R. mtr."
$ws.Cells.Item(12,3).Value = 87
$ws.Cells.Item(12,4).Value = "'17"
$ws.Cells.Item(12,5).Value = "'25 mm"
$ws.Cells.Item(12,6).Value = 56
$ws.Cells.Item(12,7).Value = "'4872.00"

# Row 13
$ws.Cells.Item(13,1).Value = "'Mtr."
$ws.Cells.Item(13,3).Value = 22
$ws.Cells.Item(13,4).Value = "'23"
$ws.Cells.Item(13,5).Value = "'8 SWG G.I. ( Hot Dipped  ) Wire "
$ws.Cells.Item(13,6).Value = 20
$ws.Cells.Item(13,7).Value = "'440.00"

# Row 14
$ws.Cells.Item(14,1).Value = "'Each"
$ws.Cells.Item(14,3).Value = 3
$ws.Cells.Item(14,4).Value = "'30"
$ws.Cells.Item(14,5).Value = "' 6 A to 32 A rating"
$ws.Cells.Item(14,6).Value = 187
$ws.Cells.Item(14,7).Value = "'561.00"

# Row 15
$ws.Cells.Item(15,1).Value = "'"
$ws.Cells.Item(15,3).Value = 63
$ws.Cells.Item(15,4).Value = "'18.0"
$ws.Cells.Item(15,5).Value = "'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = "'0.00"

# Row 16
$ws.Cells.Item(16,3).Value = 81

# Grand-total rows recompute to the new column-G sum
$ws.Cells.Item(18,7).Value = "'22149.00"
$ws.Cells.Item(18,8).Value = "'22149.00"
$ws.Cells.Item(20,7).Value = "'22149.00"
$ws.Cells.Item(20,8).Value = "'22149.00"
